# VCF loader for EAV schema
# Adds 14 new attribute rows (115-128) describing VCF-derived EAV attributes,
# extends the helper "INSERT INTO dbo.attributes ..." formula down to match,
# and leaves the selection/view pointed at the newly-added tail of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# id, Name, Value Type, Code, Code System
$rows = @(
    @(114, "Variant Call Format variant",   "binary",     $null,        $null),
    @(115, "Reference base",                "short_text", $null,        $null),
    @(116, "Number of Samples With Data",   "int",        "INFO:NS",    "VCF"),
    @(117, "Total Depth",                   "int",        "INFO:DP",    "VCF"),
    @(118, "Allele Frequency",              "float",      "INFO:AF",    "VCF"),
    @(119, "Ancestral Allele",              "short_text", "INFO:AA",    "VCF"),
    @(120, "dbSNP membership, build 129",   "binary",     "INFO:DB",    "VCF"),
    @(121, "HapMap2 membership",            "binary",     "INFO:H2",    "VCF"),
    @(122, "Genotype Quality",              "int",        "FORMAT:GQ",  "VCF"),
    @(123, "Genotype",                      "short_text", "FORMAT:GT",  "VCF"),
    @(124, "Read Depth",                    "int",        "FORMAT:DP",  "VCF"),
    @(125, "q10 Filter",                    "binary",     "FILTER:q10", "VCF"),
    @(126, "s50 Filter",                    "binary",     "FILTER:s50", "VCF"),
    @(127, "Quality",                       "float",      "Quality",    "VCF")
)

$startRow = 115
$r = $startRow
foreach ($row in $rows) {
    $id = $row[0]
    $name = $row[1]
    $valueType = $row[2]
    $code = $row[3]
    $codeSystem = $row[4]

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $valueType
    if ($code -ne $null) {
        $ws.Cells.Item($r, 4).Value = $code
    }
    if ($codeSystem -ne $null) {
        $ws.Cells.Item($r, 5).Value = $codeSystem
    }

    $formula = '=CONCATENATE("INSERT INTO dbo.attributes (id, name, value_type, code, code_system) VALUES (", A' + $r + ', ", ''", B' + $r + ', "'', ''", C' + $r + ', "'', ", IF(D' + $r + ' = "", "NULL", CONCATENATE("''", D' + $r + ', "''")), ", ", IF(E' + $r + ' = "", "NULL", CONCATENATE("''", E' + $r + ', "''")), ")")'
    $ws.Cells.Item($r, 7).Formula = $formula

    $r = $r + 1
}

$lastRow = $r - 1

$wb.Application.Calculate()

# Move the view/selection to the newly added tail, matching the author's
# final cursor position after appending the VCF attribute rows.
$ws.Range("F" + $lastRow).Select()
